$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append, matching the source "MV -datos-" refresh.
$dates = @("04-09-2021", "05-09-2021", "06-09-2021", "07-09-2021")
$data = @(
    @(2205, 2751, 15914, 2778, 4810, 9240),
    @(2205, 2751, 15914, 2778, 4810, 9240),
    @(2211, 2758, 15952, 2784, 4821, 9263),
    @(2200, 2744, 15873, 2771, 4798, 9217)
)

$startRow = 248
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Write the date label via a text formula, then paste-special it back
    # as a value. A plain .Value/.Formula string assignment of a
    # "dd-mm-yyyy"-looking literal gets auto-converted by Excel into a date
    # serial number (with a new date number-format style), which would not
    # match the plain text/shared-string cells used throughout column A.
    # Routing it through a formula result and Paste Values keeps it a
    # literal text value instead.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = '="' + $dates[$i] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $data[$i][$c]
    }
}

$excel.CutCopyMode = $false
